$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.329722881317139
$ws.Range("B1").Value = 2.632085561752319
$ws.Range("D1").Value = 1.57459557056427
$ws.Range("E1").Value = 0.9392154216766357
